$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "duplicate_image_filename" column (E) with "NA" for data rows 2-21
$ws.Range("E2:E21").Value = "NA"
